$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B3").Value = 550
$ws.Range("B7").Value = 550
$ws.Range("B12").Value = 550
$ws.Range("B16").Value = 550
$ws.Range("B19").Value = 0.6

$ws.Range("H31").Select()
